$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number. Every data row (2..304)
# had it bumped from 45190 (2023-09-21) to 45192 (2023-09-23).
$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value2
    if ($val -eq 45190) {
        $cell.Value2 = 45192
    }
}
